{"js": "// Each entry is [oldText, newText]. The document consists of one title\n// paragraph (\"YYYY-MM-DD Weekday\") followed by a 20x5 table of simple\n// arithmetic expressions (\"a+b=\" / \"a-b=\"); every paragraph in the body\n// (title + 100 table-cell paragraphs, in document order) holds exactly\n// one such run. All old values are unique, so positional replacement\n// driven by paragraph order is safe and matches the source diff 1:1.\nconst replacements = [\n  [\"2023-05-01 Monday\", \"2023-05-02 Tuesday\"],\n  [\"84+10=\", \"50+0=\"],\n  [\"56-36=\", \"4+60=\"],\n  [\"9+5=\", \"36-17=\"],\n  [\"38+3=\", \"77+18=\"],\n  [\"85-0=\", \"33+0=\"],\n  [\"50-39=\", \"70+18=\"],\n  [\"0+58=\", \"85-59=\"],\n  [\"68-11=\", \"59-29=\"],\n  [\"23+47=\", \"72-28=\"],\n  [\"44+11=\", \"10+62=\"],\n  [\"79-0=\", \"83-9=\"],\n  [\"25+73=\", \"39+50=\"],\n  [\"46-31=\", \"67+14=\"],\n  [\"27+46=\", \"76+22=\"],\n  [\"73-52=\", \"47-18=\"],\n  [\"94-60=\", \"59-45=\"],\n  [\"40-0=\", \"80-14=\"],\n  [\"2+69=\", \"81-15=\"],\n  [\"15+56=\", \"58-14=\"],\n  [\"8+56=\", \"65-16=\"],\n  [\"58-57=\", \"86-19=\"],\n  [\"9-6=\", \"40+43=\"],\n  [\"10+43=\", \"68-24=\"],\n  [\"36+25=\", \"54+13=\"],\n  [\"85+9=\", \"28+40=\"],\n  [\"71-65=\", \"1+88=\"],\n  [\"86-35=\", \"0+16=\"],\n  [\"94-57=\", \"44-35=\"],\n  [\"83-26=\", \"94+0=\"],\n  [\"70-15=\", \"3+86=\"],\n  [\"70-38=\", \"49-32=\"],\n  [\"50+14=\", \"64-51=\"],\n  [\"79-15=\", \"25+44=\"],\n  [\"97-37=\", \"71-54=\"],\n  [\"42+50=\", \"16+71=\"],\n  [\"36+8=\", \"30+25=\"],\n  [\"39+24=\", \"36-33=\"],\n  [\"13+34=\", \"98-55=\"],\n  [\"31+12=\", \"64-43=\"],\n  [\"44-9=\", \"74-11=\"],\n  [\"7+53=\", \"43-12=\"],\n  [\"68-31=\", \"47+46=\"],\n  [\"72+19=\", \"86-2=\"],\n  [\"39+31=\", \"64-63=\"],\n  [\"53-34=\", \"43+34=\"],\n  [\"34+43=\", \"55-12=\"],\n  [\"83+3=\", \"33-23=\"],\n  [\"86-70=\", \"19-19=\"],\n  [\"36-3=\", \"33+8=\"],\n  [\"45-32=\", \"74-3=\"],\n  [\"55-50=\", \"10+21=\"],\n  [\"11+57=\", \"59-49=\"],\n  [\"49-44=\", \"58-51=\"],\n  [\"45-28=\", \"87-86=\"],\n  [\"82-52=\", \"54-49=\"],\n  [\"55-46=\", \"28+47=\"],\n  [\"61-29=\", \"96-18=\"],\n  [\"11+5=\", \"43+18=\"],\n  [\"26-24=\", \"85-19=\"],\n  [\"44-27=\", \"66-19=\"],\n  [\"12+62=\", \"81-45=\"],\n  [\"71-37=\", \"59+37=\"],\n  [\"1+3=\", \"33+8=\"],\n  [\"65-27=\", \"15+29=\"],\n  [\"18+2=\", \"26+37=\"],\n  [\"59+8=\", \"82+2=\"],\n  [\"93-46=\", \"35-34=\"],\n  [\"32-4=\", \"10+71=\"],\n  [\"2+97=\", \"86-42=\"],\n  [\"62-16=\", \"2+58=\"],\n  [\"10+46=\", \"92-68=\"],\n  [\"87-57=\", \"1+75=\"],\n  [\"4+87=\", \"0+43=\"],\n  [\"79-5=\", \"75-39=\"],\n  [\"97-28=\", \"82-33=\"],\n  [\"67-17=\", \"6-1=\"],\n  [\"99-30=\", \"54+7=\"],\n  [\"82-72=\", \"46-37=\"],\n  [\"97-96=\", \"59+23=\"],\n  [\"11+60=\", \"77-12=\"],\n  [\"66-52=\", \"95-30=\"],\n  [\"68-5=\", \"91-72=\"],\n  [\"78-19=\", \"10+86=\"],\n  [\"36+46=\", \"95-36=\"],\n  [\"74-34=\", \"56-3=\"],\n  [\"91-45=\", \"61+10=\"],\n  [\"93-40=\", \"74-63=\"],\n  [\"22+43=\", \"19-14=\"],\n  [\"11-9=\", \"16+58=\"],\n  [\"91-8=\", \"20+65=\"],\n  [\"26+45=\", \"69-4=\"],\n  [\"24-7=\", \"94-1=\"],\n  [\"54+24=\", \"39+29=\"],\n  [\"80-38=\", \"5+73=\"],\n  [\"64-6=\", \"23+7=\"],\n  [\"22-19=\", \"18+76=\"],\n  [\"51+20=\", \"4+11=\"],\n  [\"25+26=\", \"62-21=\"],\n  [\"84-16=\", \"49-8=\"],\n  [\"35+61=\", \"57-30=\"]\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nif (paragraphs.items.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = paragraphs.items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Paragraph ${i}: expected \"${oldText}\" but found \"${para.text}\"`\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document body is a title paragraph (\"YYYY-MM-DD Weekday\") followed\n# by a 20x5 table of simple arithmetic expressions (\"a+b=\" / \"a-b=\"). Every\n# old value below is a unique run of text in the document, so a literal\n# Find/Replace for each pair updates exactly the one intended run.\n$replacements = @(\n    @{Old='2023-05-01 Monday'; New='2023-05-02 Tuesday'}\n    @{Old='84+10='; New='50+0='}\n    @{Old='56-36='; New='4+60='}\n    @{Old='9+5='; New='36-17='}\n    @{Old='38+3='; New='77+18='}\n    @{Old='85-0='; New='33+0='}\n    @{Old='50-39='; New='70+18='}\n    @{Old='0+58='; New='85-59='}\n    @{Old='68-11='; New='59-29='}\n    @{Old='23+47='; New='72-28='}\n    @{Old='44+11='; New='10+62='}\n    @{Old='79-0='; New='83-9='}\n    @{Old='25+73='; New='39+50='}\n    @{Old='46-31='; New='67+14='}\n    @{Old='27+46='; New='76+22='}\n    @{Old='73-52='; New='47-18='}\n    @{Old='94-60='; New='59-45='}\n    @{Old='40-0='; New='80-14='}\n    @{Old='2+69='; New='81-15='}\n    @{Old='15+56='; New='58-14='}\n    @{Old='8+56='; New='65-16='}\n    @{Old='58-57='; New='86-19='}\n    @{Old='9-6='; New='40+43='}\n    @{Old='10+43='; New='68-24='}\n    @{Old='36+25='; New='54+13='}\n    @{Old='85+9='; New='28+40='}\n    @{Old='71-65='; New='1+88='}\n    @{Old='86-35='; New='0+16='}\n    @{Old='94-57='; New='44-35='}\n    @{Old='83-26='; New='94+0='}\n    @{Old='70-15='; New='3+86='}\n    @{Old='70-38='; New='49-32='}\n    @{Old='50+14='; New='64-51='}\n    @{Old='79-15='; New='25+44='}\n    @{Old='97-37='; New='71-54='}\n    @{Old='42+50='; New='16+71='}\n    @{Old='36+8='; New='30+25='}\n    @{Old='39+24='; New='36-33='}\n    @{Old='13+34='; New='98-55='}\n    @{Old='31+12='; New='64-43='}\n    @{Old='44-9='; New='74-11='}\n    @{Old='7+53='; New='43-12='}\n    @{Old='68-31='; New='47+46='}\n    @{Old='72+19='; New='86-2='}\n    @{Old='39+31='; New='64-63='}\n    @{Old='53-34='; New='43+34='}\n    @{Old='34+43='; New='55-12='}\n    @{Old='83+3='; New='33-23='}\n    @{Old='86-70='; New='19-19='}\n    @{Old='36-3='; New='33+8='}\n    @{Old='45-32='; New='74-3='}\n    @{Old='55-50='; New='10+21='}\n    @{Old='11+57='; New='59-49='}\n    @{Old='49-44='; New='58-51='}\n    @{Old='45-28='; New='87-86='}\n    @{Old='82-52='; New='54-49='}\n    @{Old='55-46='; New='28+47='}\n    @{Old='61-29='; New='96-18='}\n    @{Old='11+5='; New='43+18='}\n    @{Old='26-24='; New='85-19='}\n    @{Old='44-27='; New='66-19='}\n    @{Old='12+62='; New='81-45='}\n    @{Old='71-37='; New='59+37='}\n    @{Old='1+3='; New='33+8='}\n    @{Old='65-27='; New='15+29='}\n    @{Old='18+2='; New='26+37='}\n    @{Old='59+8='; New='82+2='}\n    @{Old='93-46='; New='35-34='}\n    @{Old='32-4='; New='10+71='}\n    @{Old='2+97='; New='86-42='}\n    @{Old='62-16='; New='2+58='}\n    @{Old='10+46='; New='92-68='}\n    @{Old='87-57='; New='1+75='}\n    @{Old='4+87='; New='0+43='}\n    @{Old='79-5='; New='75-39='}\n    @{Old='97-28='; New='82-33='}\n    @{Old='67-17='; New='6-1='}\n    @{Old='99-30='; New='54+7='}\n    @{Old='82-72='; New='46-37='}\n    @{Old='97-96='; New='59+23='}\n    @{Old='11+60='; New='77-12='}\n    @{Old='66-52='; New='95-30='}\n    @{Old='68-5='; New='91-72='}\n    @{Old='78-19='; New='10+86='}\n    @{Old='36+46='; New='95-36='}\n    @{Old='74-34='; New='56-3='}\n    @{Old='91-45='; New='61+10='}\n    @{Old='93-40='; New='74-63='}\n    @{Old='22+43='; New='19-14='}\n    @{Old='11-9='; New='16+58='}\n    @{Old='91-8='; New='20+65='}\n    @{Old='26+45='; New='69-4='}\n    @{Old='24-7='; New='94-1='}\n    @{Old='54+24='; New='39+29='}\n    @{Old='80-38='; New='5+73='}\n    @{Old='64-6='; New='23+7='}\n    @{Old='22-19='; New='18+76='}\n    @{Old='51+20='; New='4+11='}\n    @{Old='25+26='; New='62-21='}\n    @{Old='84-16='; New='49-8='}\n    @{Old='35+61='; New='57-30='}\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $found = $find.Execute(\n        $find.Text,             # FindText\n        $false,                 # MatchCase\n        $false,                 # MatchWholeWord\n        $false,                 # MatchWildcards\n        $false,                 # MatchSoundsLike\n        $false,                 # MatchAllWordForms\n        $true,                  # Forward\n        1,                      # Wrap (wdFindContinue)\n        $false,                 # Format\n        $find.Replacement.Text, # ReplaceWith\n        2                       # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Could not find expected text: $($pair.Old)\"\n    }\n}\n"}
